# Update column C ("Förändrad") from 45726 to 45727 for rows 2 through 43
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = 43
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45726) {
        $cell.Value2 = 45727
    }
}
